$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.186.59"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.901.32"
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'307.89"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.5203"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "'0.3767"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.07275"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").Value = "'0.9042"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "'0.08328"
$ws.Range("E12").Value = "  +8.94%  "
$ws.Range("D13").Value = "1.920.47"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "'96.73"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "'0.000008653"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'0.9994"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "27.226.28"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'5.086"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "2.155.49"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "'10.64"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "'6.442"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "'146.39"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'1.752"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'18.25"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Value = "'115.11"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'4.836"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'4.895"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'0.09262"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'0.05077"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'0.7996"
$ws.Range("E34").Value = "  +4.71%  "
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "'3.422"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").Value = "'2.942"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").Value = "'2.596"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'0.5751"
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").Value = "'0.02001"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'1.080"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'9.033"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'6.606"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'117.21"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "'0.4866"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'0.9995"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'10.09"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'1.632"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'37.73"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'64.04"
$ws.Range("E51").Value = "  +0.44%  "
